$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 ("Record_map"): add a new "HeroEight" column (H) mirroring column G ---

# Copy formatting from column G into the new column H for both data blocks
$ws2.Range("G1:G13").Copy()
$ws2.Range("H1:H13").PasteSpecial(-4122)
$ws2.Range("G14:G26").Copy()
$ws2.Range("H14:H26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new header / value cells for the HeroEight column
$ws2.Range("H11").Value = "HeroEight"
$ws2.Range("H12").Value = "object"
$ws2.Range("H24").Value = "HeroEight"
$ws2.Range("H25").Value = "object"

# Swap the Row/Col sample values in both blocks
$ws2.Range("B2").Value = 7
$ws2.Range("B3").Value = 8
$ws2.Range("B15").Value = 7
$ws2.Range("B16").Value = 8

# Extend the "whole number > 0" data validation to include the new column H
$ws2.Range("C8:G9").Validation.Delete()
$ws2.Range("A2:B3").Validation.Delete()
$ws2.Range("C21:G22").Validation.Delete()
$ws2.Range("A15:B16").Validation.Delete()
$ws2.Range("C8:H9").Validation.Add(1, 1, 5, "0")
$ws2.Range("A2:B3").Validation.Add(1, 1, 5, "0")
$ws2.Range("C21:H22").Validation.Add(1, 1, 5, "0")
$ws2.Range("A15:B16").Validation.Add(1, 1, 5, "0")

# --- Window / selection changes ---

# Sheet1 loses its selected-tab marker; sheet2 becomes the active sheet/tab
# (this also drives the workbook-level activeTab value on save).
$ws2.Activate()
$ws2.Range("D31").Select()

# Best-effort: match the window geometry recorded in the commit (not every
# COM host persists these window-placement values, but setting them is
# harmless if ignored).
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
